# "pictures en cape y exchas" — add a "Fotos (IF)" column (L) to the header
# row of the DATOS sheet, with matching thin vertical separators between the
# existing header cells.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. New header cell: L1 = "Fotos (IF)" -------------------------------
$ws.Range("L1").Value = "Fotos (IF)"

# --- 2. Borders: add thin vertical separators between header cells -------
# xlEdgeLeft = 7, xlEdgeRight = 10 ; xlContinuous = 1
# A1 already has medium left/top/bottom -> just add a thin right edge.
$ws.Range("A1").Borders.Item(10).LineStyle = 1

# B1..K1 currently have no left/right border -> add thin on both sides.
$ws.Range("B1:K1").Borders.Item(7).LineStyle = 1
$ws.Range("B1:K1").Borders.Item(10).LineStyle = 1

# L1 (new last column) takes over the medium right edge; give it a thin left edge.
$ws.Range("L1").Borders.Item(7).LineStyle = 1
$ws.Range("L1").Borders.Item(10).LineStyle = 1
$ws.Range("L1").Borders.Item(10).Weight = -4138
$ws.Range("L1").Borders.Item(8).LineStyle = 1
$ws.Range("L1").Borders.Item(8).Weight = -4138
$ws.Range("L1").Borders.Item(9).LineStyle = 1
$ws.Range("L1").Borders.Item(9).Weight = -4138

# K1 no longer owns the outer-right medium edge (L1 does now) -> thin it out.
$ws.Range("K1").Borders.Item(10).LineStyle = 1

# --- 3. Alignment: match each header cell's wrap/center formatting -------
$ws.Range("L1").HorizontalAlignment = -4108
$ws.Range("L1").WrapText = $true

Write-Output "done"
